$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 22-30 with new values ---

# Row 22
$ws.Range("D22").Value = [DateTime]"2021-10-20"
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 9500
$ws.Range("P22").Value = 190

# Row 23
$ws.Range("D23").Value = [DateTime]"2021-10-20"
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8500
$ws.Range("M23").Value = 8250
$ws.Range("P23").Value = 275

# Row 24
$ws.Range("D24").Value = [DateTime]"2021-06-11"
$ws.Range("H24").Value = "Argentina(o)"
$ws.Range("K24").Value = 18000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 19000
$ws.Range("N24").Value = "$/caja 50 unidades"
$ws.Range("P24").Value = 380
$ws.Range("Q24").Value = 50

# Row 25
$ws.Range("D25").Value = [DateTime]"2021-06-11"
$ws.Range("H25").Value = "Española"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 18000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19000
$ws.Range("N25").Value = "$/caja 30 unidades"
$ws.Range("P25").Value = 633
$ws.Range("Q25").Value = 30

# Row 26
$ws.Range("D26").Value = [DateTime]"2021-05-19"
$ws.Range("J26").Value = 100

# Row 27
$ws.Range("D27").Value = [DateTime]"2021-06-03"
$ws.Range("H27").Value = "Argentina(o)"
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 15000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 15600
$ws.Range("N27").Value = "$/caja 50 unidades"
$ws.Range("P27").Value = 312
$ws.Range("Q27").Value = 50

# Row 28
$ws.Range("D28").Value = [DateTime]"2021-06-03"
$ws.Range("H28").Value = "Española"
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 17000
$ws.Range("L28").Value = 18000
$ws.Range("M28").Value = 17500
$ws.Range("N28").Value = "$/caja 30 unidades"
$ws.Range("P28").Value = 583
$ws.Range("Q28").Value = 30

# Row 29
$ws.Range("D29").Value = [DateTime]"2021-07-08"
$ws.Range("H29").Value = "Española"
$ws.Range("K29").Value = 17000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 17500
$ws.Range("N29").Value = "$/caja 30 unidades"
$ws.Range("P29").Value = 583
$ws.Range("Q29").Value = 30

# Row 30
$ws.Range("D30").Value = [DateTime]"2021-08-27"
$ws.Range("H30").Value = "Argentina(o)"
$ws.Range("N30").Value = "$/caja 50 unidades"
$ws.Range("P30").Value = 290
$ws.Range("Q30").Value = 50

# --- Append new rows 31 and 32 (copied from old row 30 data, split into two rows) ---

# Row 31 - copy formatting from row 30, then set values
$ws.Range("A30:R30").Copy()
$ws.Range("A31:R31").PasteSpecial(-4122)
$ws.Range("A31:R31").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = [DateTime]"2021-08-25"
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100112013
$ws.Range("G31").Value = "Alcachofa"
$ws.Range("H31").Value = "Argentina(o)"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 14000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 14500
$ws.Range("N31").Value = "$/caja 50 unidades"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 290
$ws.Range("Q31").Value = 50
$ws.Range("R31").Value = "Hortaliza"

# Row 32 - copy formatting from row 30, then set values
$ws.Range("A30:R30").Copy()
$ws.Range("A32:R32").PasteSpecial(-4122)
$ws.Range("A32:R32").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = [DateTime]"2021-07-20"
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112013
$ws.Range("G32").Value = "Alcachofa"
$ws.Range("H32").Value = "Española"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 14000
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = 14500
$ws.Range("N32").Value = "$/caja 30 unidades"
$ws.Range("O32").Value = "Provincia de Limarí"
$ws.Range("P32").Value = 483
$ws.Range("Q32").Value = 30
$ws.Range("R32").Value = "Hortaliza"
